$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Action caseworker" header (I1) - copy style from H1, then set text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "Action caseworker"

# Add "Attendance Notes" header (J1) - copy style from H1, then set text
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "Attendance Notes"

# Add data row values (no special style, matching H2)
$ws.Range("I2").Value = "abeaman"
$ws.Range("J2").Value = "notes"

# Move the active selection to J2 (last cell edited)
[void]$ws.Range("J2").Select()
